# UD-SRS-traceability.xlsx update
#
# Requirement "UD-SRS-15" ("The output of the system shall be graphically
# represented") is removed from the traceability table. Its associated CRS
# references (UD-CRS-06 / UD-CRS-07) are folded into the CRS# column of
# UD-SRS-08 (row 9), which now reads "UD-CRS-03 / UD-CRS-06 / UD-CRS-07".
# The row that used to hold UD-SRS-15 (row 16) is cleared out (left blank,
# keeping its original formatting) instead of the rows below shifting up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (UD-SRS-08): fold UD-CRS-06 / UD-CRS-07 into the CRS# column ---
$ws.Range("C9").Value = "UD-CRS-03`nUD-CRS-06`nUD-CRS-07"
$ws.Range("C9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 51

# --- Row 16 (was UD-SRS-15): clear the requirement, keep the row/format ---
$ws.Range("A16").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(1).RowHeight
$ws.Rows.Item(16).AutoFit()

# --- View state: scroll down a bit and leave the selection on D15 ---
$ws.Activate()
$ws.Range("D15").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
